$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to Text format first so numeric-looking strings
# (e.g. "1.007", "223.08") are not auto-converted to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '27.073.12'
$ws.Range('E2').Value = '  -0.83%  '
$ws.Range('D3').Value = '1.699.41'
$ws.Range('E3').Value = '  -0.55%  '
$ws.Range('D4').Value = '1.007'
$ws.Range('E4').Value = '  +0.32%  '
$ws.Range('D5').Value = '223.08'
$ws.Range('E5').Value = '  -0.37%  '
$ws.Range('D6').Value = '0.5238'
$ws.Range('E6').Value = '  -1.14%  '
$ws.Range('D7').Value = '1.007'
$ws.Range('E7').Value = '  +0.23%  '
$ws.Range('D8').Value = '0.06591'
$ws.Range('E8').Value = '  +0.69%  '
$ws.Range('D9').Value = '0.2607'
$ws.Range('E9').Value = '  -1.30%  '
$ws.Range('D10').Value = '20.44'
$ws.Range('E10').Value = '  -2.24%  '
$ws.Range('D11').Value = '0.07723'
$ws.Range('E11').Value = '  +1.03%  '
$ws.Range('B12').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C12').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D12').Value = '1.935.10'
$ws.Range('E12').Value = '  -0.47%  '
$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D13').Value = '4.406'
$ws.Range('E13').Value = '  -3.48%  '
$ws.Range('B14').Value = 'WrappedEther'
$ws.Range('C14').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D14').Value = '1.703.99'
$ws.Range('E14').Value = '  -0.24%  '
$ws.Range('D15').Value = '0.5707'
$ws.Range('E15').Value = '  -0.65%  '
$ws.Range('D16').Value = '0.0₅8098'
$ws.Range('E16').Value = '  -0.99%  '
$ws.Range('D17').Value = '66.54'
$ws.Range('E17').Value = '  -1.11%  '
$ws.Range('D18').Value = '27.118.56'
$ws.Range('E18').Value = '  -0.60%  '
$ws.Range('D19').Value = '217.23'
$ws.Range('E19').Value = '  +0.62%  '
$ws.Range('D20').Value = '1.007'
$ws.Range('E20').Value = '  +0.37%  '
$ws.Range('D21').Value = '4.569'
$ws.Range('E21').Value = '  -2.13%  '
$ws.Range('D22').Value = '10.27'
$ws.Range('E22').Value = '  -1.82%  '
$ws.Range('D23').Value = '5.990'
$ws.Range('E23').Value = '  +0.44%  '
$ws.Range('D24').Value = '1.008'
$ws.Range('E24').Value = '  +0.28%  '
$ws.Range('D25').Value = '144.77'
$ws.Range('E25').Value = '  +1.73%  '
$ws.Range('D26').Value = '1.728'
$ws.Range('E26').Value = '  -1.55%  '
$ws.Range('D27').Value = '0.1189'
$ws.Range('E27').Value = '  -2.28%  '
$ws.Range('D28').Value = '7.124'
$ws.Range('E28').Value = '  -1.99%  '
$ws.Range('D29').Value = '15.96'
$ws.Range('E29').Value = '  -2.13%  '
$ws.Range('D30').Value = '0.05284'
$ws.Range('E30').Value = '  -1.60%  '
$ws.Range('E31').Value = '  -0.57%  '
$ws.Range('D32').Value = '3.416'
$ws.Range('E32').Value = '  -1.93%  '
$ws.Range('D33').Value = '3.302'
$ws.Range('E33').Value = '  -3.31%  '
$ws.Range('D34').Value = '1.618'
$ws.Range('E34').Value = '  -1.21%  '
$ws.Range('D35').Value = '2.813'
$ws.Range('E35').Value = '  -2.26%  '
$ws.Range('B36').Value = 'ARBITRUM'
$ws.Range('C36').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D36').Value = '0.9431'
$ws.Range('E36').Value = '  -1.14%  '
$ws.Range('B37').Value = 'HuobiToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D37').Value = '2.395'
$ws.Range('E37').Value = '  -1.12%  '
$ws.Range('D38').Value = '0.5804'
$ws.Range('E38').Value = '  -1.01%  '
$ws.Range('D39').Value = '1.178.74'
$ws.Range('E39').Value = '  +13.02%  '
$ws.Range('D40').Value = '0.01631'
$ws.Range('E40').Value = '  +0.46%  '
$ws.Range('D41').Value = '1.007'
$ws.Range('E41').Value = '  +0.25%  '
$ws.Range('D42').Value = '5.719'
$ws.Range('E42').Value = '  -2.74%  '
$ws.Range('D43').Value = '0.8373'
$ws.Range('E43').Value = '  -0.24%  '
$ws.Range('D44').Value = '100.52'
$ws.Range('E44').Value = '  -0.50%  '
$ws.Range('D45').Value = '1.842.01'
$ws.Range('E45').Value = '  -0.53%  '
$ws.Range('D46').Value = '0.0₈111'
$ws.Range('E46').Value = '  -3.71%  '
$ws.Range('D47').Value = '56.84'
$ws.Range('E47').Value = '  -2.03%  '
$ws.Range('D48').Value = '0.4536'
$ws.Range('E48').Value = '  +0.73%  '
$ws.Range('D49').Value = '1.006'
$ws.Range('E49').Value = '  +0.36%  '
$ws.Range('D50').Value = '8.054'
$ws.Range('E50').Value = '  +0.01%  '
$ws.Range('D51').Value = '0.05214'
$ws.Range('E51').Value = '  -0.44%  '

# Restore original (default/normal) style on column D so no stray
# number-format style attribute is left on the cells.
$ws.Range("D2:D51").Style = "Normal"
